$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from C1 into D1:E1, then set header values
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4

# Update data rows 2-67 for columns B (2), C (3), D (4), E (5)
$ws.Cells.Item(2, 2).Value = -0.3260277978606521
$ws.Cells.Item(2, 3).Value = -0.3113347074289435
$ws.Cells.Item(2, 4).Value = -0.2970736233713464
$ws.Cells.Item(2, 5).Value = -0.2835723545011988
$ws.Cells.Item(3, 2).Value = 0.1974451442080811
$ws.Cells.Item(3, 3).Value = 0.2065752829286065
$ws.Cells.Item(3, 4).Value = 0.2166491866596918
$ws.Cells.Item(3, 5).Value = 0.2267225736545723
$ws.Cells.Item(4, 2).Value = 0.1305808780685744
$ws.Cells.Item(4, 3).Value = 0.1421566416881549
$ws.Cells.Item(4, 4).Value = 0.1546844081901775
$ws.Cells.Item(4, 5).Value = 0.1672166424599572
$ws.Cells.Item(5, 2).Value = -0.08640099787666576
$ws.Cells.Item(5, 3).Value = -0.07694397160759842
$ws.Cells.Item(5, 4).Value = -0.06782257657061308
$ws.Cells.Item(5, 5).Value = -0.05954051422487723
$ws.Cells.Item(6, 2).Value = 0.1271961316647943
$ws.Cells.Item(6, 3).Value = 0.1429155484052026
$ws.Cells.Item(6, 4).Value = 0.1581218431766338
$ws.Cells.Item(6, 5).Value = 0.1719715137667233
$ws.Cells.Item(7, 2).Value = -0.462073158579521
$ws.Cells.Item(7, 3).Value = -0.4527683761063582
$ws.Cells.Item(7, 4).Value = -0.4429215199565364
$ws.Cells.Item(7, 5).Value = -0.4329974376358153
$ws.Cells.Item(8, 2).Value = -0.2715765964713239
$ws.Cells.Item(8, 3).Value = -0.2598008161228147
$ws.Cells.Item(8, 4).Value = -0.2472124080800005
$ws.Cells.Item(8, 5).Value = -0.2345398796364273
$ws.Cells.Item(9, 2).Value = -0.4213748423931841
$ws.Cells.Item(9, 3).Value = -0.411365857694458
$ws.Cells.Item(9, 4).Value = -0.3990785898381406
$ws.Cells.Item(9, 5).Value = -0.3855349601999202
$ws.Cells.Item(10, 2).Value = 0.3436617211573456
$ws.Cells.Item(10, 3).Value = 0.3583351033494231
$ws.Cells.Item(10, 4).Value = 0.372151296262771
$ws.Cells.Item(10, 5).Value = 0.3844866954156452
$ws.Cells.Item(11, 2).Value = -0.2483358651284864
$ws.Cells.Item(11, 3).Value = -0.239615135361616
$ws.Cells.Item(11, 4).Value = -0.2306716941534318
$ws.Cells.Item(11, 5).Value = -0.2221088093088014
$ws.Cells.Item(12, 2).Value = -0.1094783026168914
$ws.Cells.Item(12, 3).Value = -0.0872659943667467
$ws.Cells.Item(12, 4).Value = -0.06715077251264258
$ws.Cells.Item(12, 5).Value = -0.04962863343179166
$ws.Cells.Item(13, 2).Value = -0.01157867802692844
$ws.Cells.Item(13, 3).Value = -0.01163689161266203
$ws.Cells.Item(13, 4).Value = -0.01119816786007568
$ws.Cells.Item(13, 5).Value = -0.01076816607872993
$ws.Cells.Item(14, 2).Value = 0.011889263515418
$ws.Cells.Item(14, 3).Value = 0.05201512159417407
$ws.Cells.Item(14, 4).Value = 0.08852251540343356
$ws.Cells.Item(14, 5).Value = 0.1206184983734711
$ws.Cells.Item(15, 2).Value = -0.05801708454880922
$ws.Cells.Item(15, 3).Value = -0.01353990189783741
$ws.Cells.Item(15, 4).Value = 0.02570676654085348
$ws.Cells.Item(15, 5).Value = 0.05919080618255411
$ws.Cells.Item(16, 2).Value = 0.2861528088080422
$ws.Cells.Item(16, 3).Value = 0.3474641149721687
$ws.Cells.Item(16, 4).Value = 0.4024415838960021
$ws.Cells.Item(16, 5).Value = 0.4502144255124137
$ws.Cells.Item(17, 2).Value = 0.5627694800712454
$ws.Cells.Item(17, 3).Value = 0.5981113559656567
$ws.Cells.Item(17, 4).Value = 0.6280971277039247
$ws.Cells.Item(17, 5).Value = 0.6523012256257039
$ws.Cells.Item(18, 2).Value = 0.0509146135411717
$ws.Cells.Item(18, 3).Value = 0.03848846692182706
$ws.Cells.Item(18, 4).Value = 0.02771043692032009
$ws.Cells.Item(18, 5).Value = 0.01795779568995341
$ws.Cells.Item(19, 2).Value = 0.3801368142500762
$ws.Cells.Item(19, 3).Value = 0.3958525534261808
$ws.Cells.Item(19, 4).Value = 0.409797240493623
$ws.Cells.Item(19, 5).Value = 0.4212192338181458
$ws.Cells.Item(20, 2).Value = 0.2357545280472422
$ws.Cells.Item(20, 3).Value = 0.2938656074185948
$ws.Cells.Item(20, 4).Value = 0.3450655875566412
$ws.Cells.Item(20, 5).Value = 0.3887353086973739
$ws.Cells.Item(21, 2).Value = 0.4612034855031213
$ws.Cells.Item(21, 3).Value = 0.5182927104842734
$ws.Cells.Item(21, 4).Value = 0.5673456989925195
$ws.Cells.Item(21, 5).Value = 0.6080725634857478
$ws.Cells.Item(22, 2).Value = 0.3407779426971171
$ws.Cells.Item(22, 3).Value = 0.3774294530484647
$ws.Cells.Item(22, 4).Value = 0.4089951628219275
$ws.Cells.Item(22, 5).Value = 0.4349674922936272
$ws.Cells.Item(23, 2).Value = -0.09213597395040349
$ws.Cells.Item(23, 3).Value = -0.06047098978932161
$ws.Cells.Item(23, 4).Value = -0.03356228790576768
$ws.Cells.Item(23, 5).Value = -0.01171677136219079
$ws.Cells.Item(24, 2).Value = 4.536142022674653
$ws.Cells.Item(24, 3).Value = 4.575953720226746
$ws.Cells.Item(24, 4).Value = 4.533580557498998
$ws.Cells.Item(24, 5).Value = 4.424748944720001
$ws.Cells.Item(25, 2).Value = 0.5264117600795265
$ws.Cells.Item(25, 3).Value = 0.4916963026454615
$ws.Cells.Item(25, 4).Value = 0.4602768751885618
$ws.Cells.Item(25, 5).Value = 0.4321816775789644
$ws.Cells.Item(26, 2).Value = 0.3854558869822713
$ws.Cells.Item(26, 3).Value = 0.3637634469251554
$ws.Cells.Item(26, 4).Value = 0.3437026357892849
$ws.Cells.Item(26, 5).Value = 0.322793265433509
$ws.Cells.Item(27, 2).Value = 0.3323569593629884
$ws.Cells.Item(27, 3).Value = 0.3013268639216591
$ws.Cells.Item(27, 4).Value = 0.2739135470728835
$ws.Cells.Item(27, 5).Value = 0.2467306723831639
$ws.Cells.Item(28, 2).Value = 1.099255200858447
$ws.Cells.Item(28, 3).Value = 1.068469446076402
$ws.Cells.Item(28, 4).Value = 1.037001730912954
$ws.Cells.Item(28, 5).Value = 1.005564956752435
$ws.Cells.Item(29, 2).Value = 5.868063101441583
$ws.Cells.Item(29, 3).Value = 5.453715541498916
$ws.Cells.Item(29, 4).Value = 5.047154797343312
$ws.Cells.Item(29, 5).Value = 4.652447346813815
$ws.Cells.Item(30, 2).Value = 1.019959663470545
$ws.Cells.Item(30, 3).Value = 0.9684969965365088
$ws.Cells.Item(30, 4).Value = 0.9206352240028133
$ws.Cells.Item(30, 5).Value = 0.8765069721295897
$ws.Cells.Item(31, 2).Value = -0.1327971309754888
$ws.Cells.Item(31, 3).Value = -0.19829785748661
$ws.Cells.Item(31, 4).Value = -0.2528510393493318
$ws.Cells.Item(31, 5).Value = -0.3009271808733578
$ws.Cells.Item(32, 2).Value = 0.8129124910470404
$ws.Cells.Item(32, 3).Value = 0.7771228137702961
$ws.Cells.Item(32, 4).Value = 0.744430169787492
$ws.Cells.Item(32, 5).Value = 0.7130909137234249
$ws.Cells.Item(33, 2).Value = 0.9408763627916946
$ws.Cells.Item(33, 3).Value = 0.9117538597084132
$ws.Cells.Item(33, 4).Value = 0.8858350655748324
$ws.Cells.Item(33, 5).Value = 0.8619185721803341
$ws.Cells.Item(34, 2).Value = -0.6237788977627812
$ws.Cells.Item(34, 3).Value = -0.6548704237507494
$ws.Cells.Item(34, 4).Value = -0.6836312085549999
$ws.Cells.Item(34, 5).Value = -0.7102928388237109
$ws.Cells.Item(35, 2).Value = 0.8227024093626152
$ws.Cells.Item(35, 3).Value = 0.813604258629973
$ws.Cells.Item(35, 4).Value = 0.8056797875994701
$ws.Cells.Item(35, 5).Value = 0.798334102803591
$ws.Cells.Item(36, 2).Value = 0.7852006021171384
$ws.Cells.Item(36, 3).Value = 0.7675934109401431
$ws.Cells.Item(36, 4).Value = 0.7526232544134352
$ws.Cells.Item(36, 5).Value = 0.7396484170207569
$ws.Cells.Item(37, 2).Value = 0.7616093557509258
$ws.Cells.Item(37, 3).Value = 0.7417926700373553
$ws.Cells.Item(37, 4).Value = 0.7248622428537608
$ws.Cells.Item(37, 5).Value = 0.7101540691582918
$ws.Cells.Item(38, 2).Value = 0.7364618531331382
$ws.Cells.Item(38, 3).Value = 0.7173287459599338
$ws.Cells.Item(38, 4).Value = 0.700641982288517
$ws.Cells.Item(38, 5).Value = 0.6857196939546594
$ws.Cells.Item(39, 2).Value = 0.5854033188892174
$ws.Cells.Item(39, 3).Value = 0.5810806496389941
$ws.Cells.Item(39, 4).Value = 0.5777955608229628
$ws.Cells.Item(39, 5).Value = 0.5749546883777493
$ws.Cells.Item(40, 2).Value = 0.7554000423666646
$ws.Cells.Item(40, 3).Value = 0.7522695508355323
$ws.Cells.Item(40, 4).Value = 0.7495143856711403
$ws.Cells.Item(40, 5).Value = 0.7465678738711784
$ws.Cells.Item(41, 2).Value = 0.5625794058477487
$ws.Cells.Item(41, 3).Value = 0.5544369560187997
$ws.Cells.Item(41, 4).Value = 0.5485379151631973
$ws.Cells.Item(41, 5).Value = 0.5440696365417749
$ws.Cells.Item(42, 2).Value = 0.7176215966072691
$ws.Cells.Item(42, 3).Value = 0.690163836811005
$ws.Cells.Item(42, 4).Value = 0.6659896593558211
$ws.Cells.Item(42, 5).Value = 0.6445022733954443
$ws.Cells.Item(43, 2).Value = 0.7239792128051638
$ws.Cells.Item(43, 3).Value = 0.7090153754814805
$ws.Cells.Item(43, 4).Value = 0.6964002319199823
$ws.Cells.Item(43, 5).Value = 0.6853986136278564
$ws.Cells.Item(44, 2).Value = 0.6829484627113083
$ws.Cells.Item(44, 3).Value = 0.6744744297863949
$ws.Cells.Item(44, 4).Value = 0.6679388035570271
$ws.Cells.Item(44, 5).Value = 0.6626134510921413
$ws.Cells.Item(45, 2).Value = 0.6771656920013283
$ws.Cells.Item(45, 3).Value = 0.6581781264025099
$ws.Cells.Item(45, 4).Value = 0.6423658368988955
$ws.Cells.Item(45, 5).Value = 0.629119162747515
$ws.Cells.Item(46, 2).Value = -1.25897553861914
$ws.Cells.Item(46, 3).Value = -1.263522051325703
$ws.Cells.Item(46, 4).Value = -1.267422265526079
$ws.Cells.Item(46, 5).Value = -1.270612138593214
$ws.Cells.Item(47, 2).Value = -0.9738671067461906
$ws.Cells.Item(47, 3).Value = -0.9793198897550003
$ws.Cells.Item(47, 4).Value = -0.9839642426851458
$ws.Cells.Item(47, 5).Value = -0.987847171869874
$ws.Cells.Item(48, 2).Value = -0.8666450125273316
$ws.Cells.Item(48, 3).Value = -0.8733249578619043
$ws.Cells.Item(48, 4).Value = -0.8783409952576358
$ws.Cells.Item(48, 5).Value = -0.8819130794400558
$ws.Cells.Item(49, 2).Value = -0.6384770605229276
$ws.Cells.Item(49, 3).Value = -0.6423388394497972
$ws.Cells.Item(49, 4).Value = -0.6449706923609597
$ws.Cells.Item(49, 5).Value = -0.6466568028909876
$ws.Cells.Item(50, 2).Value = -0.04719667975478657
$ws.Cells.Item(50, 3).Value = -0.04870843203017888
$ws.Cells.Item(50, 4).Value = -0.04975729090437483
$ws.Cells.Item(50, 5).Value = -0.05075605151621303
$ws.Cells.Item(51, 2).Value = -0.8594604294715984
$ws.Cells.Item(51, 3).Value = -0.8651112230871479
$ws.Cells.Item(51, 4).Value = -0.8694247176179918
$ws.Cells.Item(51, 5).Value = -0.8725182445532166
$ws.Cells.Item(52, 2).Value = -0.8594604294715984
$ws.Cells.Item(52, 3).Value = -0.8651112230871479
$ws.Cells.Item(52, 4).Value = -0.8694247176179918
$ws.Cells.Item(52, 5).Value = -0.8725182445532166
$ws.Cells.Item(53, 2).Value = -1.083426097556684
$ws.Cells.Item(53, 3).Value = -1.097167705592106
$ws.Cells.Item(53, 4).Value = -1.108488237167732
$ws.Cells.Item(53, 5).Value = -1.117547215306057
$ws.Cells.Item(54, 2).Value = -0.1872636393394086
$ws.Cells.Item(54, 3).Value = -0.1869482209713249
$ws.Cells.Item(54, 4).Value = -0.1860485087245954
$ws.Cells.Item(54, 5).Value = -0.1849424187574051
$ws.Cells.Item(55, 2).Value = -0.9945094100735348
$ws.Cells.Item(55, 3).Value = -0.9991769555261637
$ws.Cells.Item(55, 4).Value = -1.003171825398653
$ws.Cells.Item(55, 5).Value = -1.006553059405665
$ws.Cells.Item(56, 2).Value = -0.887447096790594
$ws.Cells.Item(56, 3).Value = -0.8841450427105935
$ws.Cells.Item(56, 4).Value = -0.8817790896142385
$ws.Cells.Item(56, 5).Value = -0.8802859356151449
$ws.Cells.Item(57, 2).Value = -0.9465834676324881
$ws.Cells.Item(57, 3).Value = -0.9339108452677687
$ws.Cells.Item(57, 4).Value = -0.9226930628719456
$ws.Cells.Item(57, 5).Value = -0.9129627576960649
$ws.Cells.Item(58, 2).Value = -1.150642723629216
$ws.Cells.Item(58, 3).Value = -1.128012501669319
$ws.Cells.Item(58, 4).Value = -1.108173379963219
$ws.Cells.Item(58, 5).Value = -1.090976065439593
$ws.Cells.Item(59, 2).Value = -0.8595943936850527
$ws.Cells.Item(59, 3).Value = -0.8479030779331405
$ws.Cells.Item(59, 4).Value = -0.8369344785974706
$ws.Cells.Item(59, 5).Value = -0.8267199931637218
$ws.Cells.Item(60, 2).Value = -0.5178421633847863
$ws.Cells.Item(60, 3).Value = -0.4977996238259934
$ws.Cells.Item(60, 4).Value = -0.4795569449529267
$ws.Cells.Item(60, 5).Value = -0.4632992899657408
$ws.Cells.Item(61, 2).Value = 0.3694990498946619
$ws.Cells.Item(61, 3).Value = 0.3714301207007228
$ws.Cells.Item(61, 4).Value = 0.3739035917335791
$ws.Cells.Item(61, 5).Value = 0.37627040402272
$ws.Cells.Item(62, 2).Value = -1.228845839492119
$ws.Cells.Item(62, 3).Value = -1.209516136962105
$ws.Cells.Item(62, 4).Value = -1.192543225469806
$ws.Cells.Item(62, 5).Value = -1.177900817433407
$ws.Cells.Item(63, 2).Value = -0.7647483462113032
$ws.Cells.Item(63, 3).Value = -0.7390840125550109
$ws.Cells.Item(63, 4).Value = -0.7147337912404659
$ws.Cells.Item(63, 5).Value = -0.6918662570581363
$ws.Cells.Item(64, 2).Value = -0.9085383082165912
$ws.Cells.Item(64, 3).Value = -0.9025893591691596
$ws.Cells.Item(64, 4).Value = -0.8955302735399159
$ws.Cells.Item(64, 5).Value = -0.8877613539810557
$ws.Cells.Item(65, 2).Value = -0.1307796336478854
$ws.Cells.Item(65, 3).Value = -0.1116670023084185
$ws.Cells.Item(65, 4).Value = -0.09428881291620919
$ws.Cells.Item(65, 5).Value = -0.07896190354401719
$ws.Cells.Item(66, 2).Value = -0.8143641485938967
$ws.Cells.Item(66, 3).Value = -0.7907523108579071
$ws.Cells.Item(66, 4).Value = -0.7706267216278589
$ws.Cells.Item(66, 5).Value = -0.753939938861058
$ws.Cells.Item(67, 2).Value = -0.7944889528103205
$ws.Cells.Item(67, 3).Value = -0.7623398618387645
$ws.Cells.Item(67, 4).Value = -0.7358762115270424
$ws.Cells.Item(67, 5).Value = -0.714736881995928
